$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-03 Monday" "2025-03-04 Tuesday"

Replace-Text "604÷4=151, 0" "577÷2=288, 1"
Replace-Text "814÷2=407, 0" "531÷3=177, 0"
Replace-Text "904÷9=100, 4" "595÷7=85, 0"
Replace-Text "930÷8=116, 2" "691÷5=138, 1"
Replace-Text "714÷9=79, 3" "866÷7=123, 5"

Replace-Text "256÷5=51, 1" "842÷8=105, 2"
Replace-Text "583÷6=97, 1" "950÷5=190, 0"
Replace-Text "745÷9=82, 7" "674÷7=96, 2"
Replace-Text "187÷8=23, 3" "985÷3=328, 1"
Replace-Text "105÷8=13, 1" "483÷9=53, 6"

Replace-Text "983÷4=245, 3" "336÷5=67, 1"
Replace-Text "503÷4=125, 3" "501÷4=125, 1"
Replace-Text "536÷2=268, 0" "220÷5=44, 0"
Replace-Text "856÷6=142, 4" "252÷5=50, 2"
Replace-Text "516÷9=57, 3" "290÷4=72, 2"

Replace-Text "470÷5=94, 0" "863÷8=107, 7"
Replace-Text "879÷7=125, 4" "164÷5=32, 4"
Replace-Text "710÷5=142, 0" "941÷5=188, 1"
Replace-Text "555÷2=277, 1" "133÷2=66, 1"
Replace-Text "243÷8=30, 3" "820÷6=136, 4"

Replace-Text "535÷6=89, 1" "620÷8=77, 4"
Replace-Text "502÷6=83, 4" "690÷9=76, 6"
Replace-Text "505÷2=252, 1" "119÷6=19, 5"
Replace-Text "512÷9=56, 8" "449÷8=56, 1"
Replace-Text "843÷4=210, 3" "766÷4=191, 2"
